$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '97.083.97'
Set-TextValue $ws.Range("E2") '  +3.61%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.323.49'
Set-TextValue $ws.Range("E3") '  +7.65%  '

# Row 4
Set-TextValue $ws.Range("E4") '  -0.09%  '

# Row 5
Set-TextValue $ws.Range("D5") '249.52'
Set-TextValue $ws.Range("E5") '  +6.69%  '

# Row 6
Set-TextValue $ws.Range("D6") '621.36'
Set-TextValue $ws.Range("E6") '  +1.86%  '

# Row 7
Set-TextValue $ws.Range("D7") '1.11'
Set-TextValue $ws.Range("E7") '  +1.37%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.384'
Set-TextValue $ws.Range("E8") '  +1.27%  '

# Row 9
Set-TextValue $ws.Range("E9") '  +0.01%  '

# Row 10
Set-TextValue $ws.Range("D10") '3.320.66'
Set-TextValue $ws.Range("E10") '  +7.62%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.785'
Set-TextValue $ws.Range("E11") '  -4.36%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.199'
Set-TextValue $ws.Range("E12") '  +1.62%  '

# Row 13
Set-TextValue $ws.Range("D13") '96.780.07'
Set-TextValue $ws.Range("E13") '  +3.20%  '

# Row 14
Set-TextValue $ws.Range("D14") '0.0000246'
Set-TextValue $ws.Range("E14") '  +2.80%  '

# Row 15
Set-TextValue $ws.Range("D15") '35.36'
Set-TextValue $ws.Range("E15") '  +3.82%  '

# Row 16
Set-TextValue $ws.Range("D16") '3.916.63'
Set-TextValue $ws.Range("E16") '  +6.76%  '

# Row 17
Set-TextValue $ws.Range("D17") '5.52'
Set-TextValue $ws.Range("E17") '  +5.35%  '

# Row 18
Set-TextValue $ws.Range("D18") '3.309.17'
Set-TextValue $ws.Range("E18") '  +5.79%  '

# Row 19
Set-TextValue $ws.Range("E19") '  -2.46%  '

# Row 20
Set-TextValue $ws.Range("D20") '14.96'
Set-TextValue $ws.Range("E20") '  +2.23%  '

# Row 21
Set-TextValue $ws.Range("D21") '482.04'
Set-TextValue $ws.Range("E21") '  +9.30%  '

# Row 22
Set-TextValue $ws.Range("B22") 'PEPE'
Set-TextValue $ws.Range("C22") 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D22") '0.0000209'
Set-TextValue $ws.Range("E22") '  +8.18%  '

# Row 23
Set-TextValue $ws.Range("B23") 'Polkadot'
Set-TextValue $ws.Range("C23") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D23") '5.82'
Set-TextValue $ws.Range("E23") '  +0.64%  '

# Row 24
Set-TextValue $ws.Range("D24") '9.26'
Set-TextValue $ws.Range("E24") '  +4.63%  '

# Row 25
Set-TextValue $ws.Range("D25") '5.66'
Set-TextValue $ws.Range("E25") '  +2.47%  '

# Row 26
Set-TextValue $ws.Range("D26") '87.70'
Set-TextValue $ws.Range("E26") '  +3.32%  '

# Row 27
Set-TextValue $ws.Range("D27") '12.10'
Set-TextValue $ws.Range("E27") '  +1.04%  '

# Row 28
Set-TextValue $ws.Range("D28") '3.490.94'
Set-TextValue $ws.Range("E28") '  +6.76%  '

# Row 29
Set-TextValue $ws.Range("E29") '  +0.08%  '

# Row 30
Set-TextValue $ws.Range("D30") '0.183'
Set-TextValue $ws.Range("E30") '  +2.06%  '

# Row 31
Set-TextValue $ws.Range("D31") '0.241'
Set-TextValue $ws.Range("E31") '  -4.57%  '

# Row 32
Set-TextValue $ws.Range("E32") '  -0.10%  '

# Row 33
Set-TextValue $ws.Range("E33") '  -0.32%  '

# Row 34
Set-TextValue $ws.Range("E34") '  -0.25%  '

# Row 35
Set-TextValue $ws.Range("D35") '27.22'
Set-TextValue $ws.Range("E35") '  +6.48%  '

# Row 36
Set-TextValue $ws.Range("D36") '7.42'
Set-TextValue $ws.Range("E36") '  -4.74%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.151'
Set-TextValue $ws.Range("E37") '  -4.76%  '

# Row 38
Set-TextValue $ws.Range("D38") '507.83'
Set-TextValue $ws.Range("E38") '  +9.14%  '

# Row 39
Set-TextValue $ws.Range("D39") '1.94'
Set-TextValue $ws.Range("E39") '  +2.77%  '

# Row 40
Set-TextValue $ws.Range("D40") '24.77'
Set-TextValue $ws.Range("E40") '  +3.53%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.452'
Set-TextValue $ws.Range("E41") '  +1.29%  '

# Row 42
Set-TextValue $ws.Range("D42") '1.29'
Set-TextValue $ws.Range("E42") '  +0.88%  '

# Row 43
Set-TextValue $ws.Range("B43") 'ARBITRUM'
Set-TextValue $ws.Range("C43") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D43") '0.800'
Set-TextValue $ws.Range("E43") '  +17.69%  '

# Row 44
Set-TextValue $ws.Range("B44") 'dogwifhat'
Set-TextValue $ws.Range("C44") 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D44") '3.26'
Set-TextValue $ws.Range("E44") '  +4.99%  '

# Row 45
Set-TextValue $ws.Range("D45") '3.51'
Set-TextValue $ws.Range("E45") '  -5.63%  '

# Row 46
Set-TextValue $ws.Range("E46") '  +0.01%  '

# Row 47
Set-TextValue $ws.Range("D47") '161.28'
Set-TextValue $ws.Range("E47") '  -0.06%  '

# Row 48
Set-TextValue $ws.Range("D48") '1.92'
Set-TextValue $ws.Range("E48") '  +3.74%  '

# Row 49
Set-TextValue $ws.Range("D49") '1.38'
Set-TextValue $ws.Range("E49") '  +7.11%  '

# Row 50
Set-TextValue $ws.Range("B50") 'Optimism'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
Set-TextValue $ws.Range("D50") '2.18'
Set-TextValue $ws.Range("E50") '  +27.87%  '

# Row 51
Set-TextValue $ws.Range("B51") 'OKB'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D51") '45.39'
Set-TextValue $ws.Range("E51") '  +3.77%  '
